# Update the daily figure in T2 and move the active selection to T2
# (matches the source workbook re-upload: value bump 475587 -> 477623,
# and the saved cursor position moving from T3 to T2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 477623
$ws.Range("T2").Select()
